# Update crypto price/volume data per the Wed Jul 5 11:11:42 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.544.43"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.920.07"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.10"
$ws.Range("E5").Value = "  -2.48%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4791"
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2884"
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06713"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.76"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "104.14"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").Value = "1.924.09"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07734"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.265"
$ws.Range("E14").Value = "  -3.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6845"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "266.57"
$ws.Range("E16").Value = "  -6.06%  "
$ws.Range("D17").Value = "30.590.19"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007533"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.77"
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.467"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.355"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.655"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.54"
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.10"
$ws.Range("E26").Value = "  -4.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.095"
$ws.Range("E27").Value = "  -5.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1025"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.389"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.534"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.520"
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.260"
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04767"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7421"
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01955"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.642"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.387"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.01"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.006"
$ws.Range("E42").Value = "  -4.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8640"
$ws.Range("E43").Value = "  -2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.50"
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4308"
$ws.Range("E45").Value = "  -3.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.557"
$ws.Range("E47").Value = "  -7.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "976.31"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1205"
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.29"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.956"
$ws.Range("E51").Value = "  -4.42%  "
